# Add two new named ranges ("i" and "j") pointing at calculations1!C10 and
# calculations1!C11, then populate a small block of cells (rows 8-16) on the
# "calculations1" sheet that exercise formulas referencing plain cell
# addresses as well as the newly (and previously) defined names - showing
# that a formula can read a cell even when its name wasn't defined ahead of
# time. Finally leave the sheet selection on C17, matching where the user
# ended up after typing this block in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calculations1")

# New defined names - pass Range objects (not plain strings) to RefersTo so
# two-digit row numbers are preserved correctly (e.g. "C10", not "10").
$wb.Names.Add("i", $ws.Range("C10"))
$wb.Names.Add("j", $ws.Range("C11"))

# New data block.
$ws.Range("C8").Value = 10
$ws.Range("C9").Formula = "=C8+1"
$ws.Range("C10").Formula = "=C8+1"
$ws.Range("C11").Value = 10

$ws.Range("C13").Formula = "=C9+1"
$ws.Range("C14").Formula = "=C10+1"
$ws.Range("C15").Formula = "=i+1"
$ws.Range("C16").Formula = "=C8+1"

# Leave the selection where the user's cursor ended up.
$ws.Activate()
$ws.Range("C17").Select() | Out-Null
